# Adds a new "priority" test-case row (row 8) to the Gila Breath test-cases
# sheet, mirroring the formatting of the last existing row (row 7), then
# updates the sheet's row-1 height and active selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 7 (values + styles) into a new row 8 by copy/insert, so the
# new row inherits the same cell styles (border/fill/font xf) already used
# by the other data rows instead of picking up the workbook default style.
$ws.Rows.Item(7).Copy()
$ws.Rows.Item(8).Insert()

# Overwrite the duplicated values with the new "priority" test case content.
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "priority"
$ws.Range("C8").Value = "checking the priority functionality"
$ws.Range("D8").Value = "1. click on priorities option in the 1st date slot tab"
$ws.Range("E8").Value = "unwanted pop up"
$ws.Range("F8").Value = "priorities tab should work"
$ws.Range("G8").Value = "pop up comes up showing the server is not working"

# Match the explicit row height used by the other data rows.
$ws.Rows.Item(8).RowHeight = 15.75

# Header row gets an explicit height too.
$ws.Rows.Item(1).RowHeight = 15

# Move the active selection to D10, as in the edited workbook.
$null = $ws.Range("D10").Select()
